# "after haunting mars #3" -- add the new "Durango" habitat to the Mars
# habs population sheet, in its correct sorted (population-descending)
# position between "Cipango" (1700) and "Pilsener City" (1000), i.e. row 37,
# pushing the previous rows 37-41 down to 38-42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room for the new row right where it belongs in the
#    population-sorted list (this also carries the little "Senate
#    composition" side table in M/N down by one row, along with every
#    other row below it).
$ws.Rows.Item(37).Insert()

# 2. Fill in the new habitat's data (same order the author typed it in:
#    Hab, Source, Region, then back to fill in Major Influencers, Location,
#    Type and the Comment).
$ws.Range("A37").Value = "Durango"
$ws.Range("D37").Value = "Haunting Mars"
$ws.Range("E37").Value = "HEL 11"
$ws.Range("C37").Value = "Fa Jing (formerly TTO)"
$ws.Range("B37").Value = "Hellas Plaintia"
$ws.Range("F37").Value = "Bubble Town"
$ws.Range("G37").Value = 1500
$ws.Range("H37").Value = 0.4
$ws.Range("I37").Formula = "=G37*H37"
$ws.Range("J37").Value = "Taken over by Faa Jing"

# 3. The row-insert above pulled the summary cells outside the main table
#    into its reference-adjustment logic even though their author never
#    touched them by hand -- restore the small "Senate composition" helper
#    formulas (columns L/N) to the exact text the sheet had before the
#    insert. Only the "total" line (which really did move down one row
#    with the table it summarizes) gets a deliberately updated range.
$ws.Range("L5").Formula = "=SUM(G1:G124)"
$ws.Range("L9").Formula = "=SUM(I3:I125)"
$ws.Range("N30").Formula = '=COUNTIF(D13:D1012,"*D*")+COUNTIF(D13:D1012,"DD*")+COUNTIF(D13:D1012,"DDD")'
$ws.Range("N31").Formula = '=COUNTIF(D13:D112,"*X*")+COUNTIF(D13:D112,"XX*")+COUNTIF(D13:D1012,"XXX")'
$ws.Range("N33").Formula = '=COUNTIF(D13:D112,"*F*")+COUNTIF(D13:D112,"FF*")+COUNTIF(D13:D1012,"FFF")'
$ws.Range("N34").Formula = '=COUNTIF(D13:D112,"*R*")+COUNTIF(D13:D112,"RR*")+COUNTIF(D13:D1012,"RRR")'
$ws.Range("N38").Formula = "=SUM(N31:N37)"

# 4. Selection / view bookkeeping to match where the author ended up.
$ws.Range("H37").Select()
